$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (rows 2-9) ---
# Row 5: "DEPARTAMENTO" -> "SISTEMAS COMPUTACIONALES"; clear other cells in the row
$ws.Range("C5:G5").Clear()
$ws.Range("E5").Value = "SISTEMAS COMPUTACIONALES"

# Row 6: keep "AÑO 2024" but clear surrounding cells
$ws.Range("C6:G6").Clear()
$ws.Range("E6").Value = "AÑO 2024"

# Row 8: "PERIODO" -> "AGOSTO - DICIEMBRE 2024"
$ws.Range("E8").Value = "AGOSTO - DICIEMBRE 2024"

# Row 9: clear the instructions text, keep formatting
$ws.Range("E9").ClearContents()

# --- Table rows ---
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 2
$ws.Range("B13").Value = 3

# Row 14 totals: replace formulas with static values
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1

# --- New signature rows 18-19 ---
$ws.Range("E18").Value = "MARITZA FLORES SARABIA"
$ws.Range("F18").Value = ""
$ws.Range("E19").Value = "JEFE DEL DEPARTAMENTO DE DESARROLLO ACADEMICO"
$ws.Range("F19").Value = ""
$ws.Range("E18:F18").Merge()
$ws.Range("E19:F19").Merge()

# --- Column widths ---
$ws.Columns("F").ColumnWidth = 31.25
$ws.Columns("E").ColumnWidth = 31.25

# --- Selection ---
$ws.Range("D5").Select()
